{"js": "// Replace the two-digit multiplication problems in the document's table\n// with the new values, as described by the commit diff.\nconst replacements = [\n  [\"38\u00d781=\", \"26\u00d755=\"],\n  [\"60\u00d775=\", \"29\u00d740=\"],\n  [\"12\u00d734=\", \"14\u00d717=\"],\n  [\"90\u00d726=\", \"20\u00d795=\"],\n  [\"60\u00d744=\", \"94\u00d756=\"],\n  [\"92\u00d773=\", \"70\u00d757=\"],\n  [\"52\u00d713=\", \"53\u00d782=\"],\n  [\"30\u00d755=\", \"76\u00d750=\"],\n  [\"74\u00d752=\", \"39\u00d768=\"],\n  [\"53\u00d795=\", \"32\u00d737=\"],\n  [\"51\u00d752=\", \"95\u00d778=\"],\n  [\"44\u00d732=\", \"30\u00d797=\"],\n  [\"57\u00d711=\", \"14\u00d722=\"],\n  [\"99\u00d740=\", \"85\u00d745=\"],\n  [\"46\u00d754=\", \"92\u00d740=\"],\n  [\"66\u00d733=\", \"44\u00d791=\"],\n  [\"26\u00d775=\", \"30\u00d716=\"],\n  [\"38\u00d772=\", \"32\u00d792=\"],\n  [\"77\u00d725=\", \"99\u00d760=\"],\n  [\"83\u00d763=\", \"34\u00d723=\"],\n  [\"78\u00d714=\", \"70\u00d792=\"],\n  [\"40\u00d785=\", \"65\u00d737=\"],\n  [\"84\u00d718=\", \"23\u00d714=\"],\n  [\"75\u00d741=\", \"43\u00d767=\"],\n  [\"41\u00d777=\", \"41\u00d769=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the document's table\n# with the new values, as described by the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"38\u00d781=\", \"26\u00d755=\"),\n  @(\"60\u00d775=\", \"29\u00d740=\"),\n  @(\"12\u00d734=\", \"14\u00d717=\"),\n  @(\"90\u00d726=\", \"20\u00d795=\"),\n  @(\"60\u00d744=\", \"94\u00d756=\"),\n  @(\"92\u00d773=\", \"70\u00d757=\"),\n  @(\"52\u00d713=\", \"53\u00d782=\"),\n  @(\"30\u00d755=\", \"76\u00d750=\"),\n  @(\"74\u00d752=\", \"39\u00d768=\"),\n  @(\"53\u00d795=\", \"32\u00d737=\"),\n  @(\"51\u00d752=\", \"95\u00d778=\"),\n  @(\"44\u00d732=\", \"30\u00d797=\"),\n  @(\"57\u00d711=\", \"14\u00d722=\"),\n  @(\"99\u00d740=\", \"85\u00d745=\"),\n  @(\"46\u00d754=\", \"92\u00d740=\"),\n  @(\"66\u00d733=\", \"44\u00d791=\"),\n  @(\"26\u00d775=\", \"30\u00d716=\"),\n  @(\"38\u00d772=\", \"32\u00d792=\"),\n  @(\"77\u00d725=\", \"99\u00d760=\"),\n  @(\"83\u00d763=\", \"34\u00d723=\"),\n  @(\"78\u00d714=\", \"70\u00d792=\"),\n  @(\"40\u00d785=\", \"65\u00d737=\"),\n  @(\"84\u00d718=\", \"23\u00d714=\"),\n  @(\"75\u00d741=\", \"43\u00d767=\"),\n  @(\"41\u00d777=\", \"41\u00d769=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
